$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, shifting rows 10-65 down to 11-66
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with data
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10, 3).Value = "Los Lagos"
$ws.Cells.Item(10, 4).Value = 44749
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 100112043
$ws.Cells.Item(10, 7).Value = "Pepino dulce"
$ws.Cells.Item(10, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 9).Value = "Especial"
$ws.Cells.Item(10, 10).Value = 35
$ws.Cells.Item(10, 11).Value = 21000
$ws.Cells.Item(10, 12).Value = 21000
$ws.Cells.Item(10, 13).Value = 21000
$ws.Cells.Item(10, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 1167
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = "Hortaliza"
